$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B58 was an inline string "3"; change it to a real number 3
$ws.Range("B58").Value = 3

# Add new row 59 with data
$ws.Range("A59").Value = "Ruilin"
$ws.Range("B59").Value = "'2"
$ws.Range("C59").Value = "The paper's contributions seem to be marginal."
$ws.Range("D59").Value = "CRT"
$ws.Range("E59").Value = "OTH"
$ws.Range("F59").Value = "0752b86f-1f11-4877-863b-8b76e37d0f73"
$ws.Range("G59").Value = "rJrTwxbCb_annotated.xlsx"
$ws.Range("H59").Value = "The paper's contributions seem to be marginal."
